# ---------------------------------------------------------------------------
# Applies the "Completed inspect raw ... main_data_io.py ... time_helpers.py"
# commit to the gvp_project_structure_and_documentation workbook.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsColumnInfo = $wb.Worksheets.Item("column_info")
$wsPyFiles    = $wb.Worksheets.Item("py_files_doc")
$wsChangelog  = $wb.Worksheets.Item("changelog")

# ---------------------------------------------------------------------------
# 1) changelog sheet: append a new row describing main_data_io.py
# ---------------------------------------------------------------------------

$wsChangelog.Range("A3:D3").Copy()
$wsChangelog.Range("A4:D4").PasteSpecial(-4122) # xlPasteFormats

$wsChangelog.Range("A4").Value = 3
$wsChangelog.Range("B4").Value = 45881
$wsChangelog.Range("C4").Value = "main_data_io.py"
$wsChangelog.Range("D4").Value = "1. Incorporated loading of raw data for initial inspection and standardization.`n2. Incorporated saving of cleaned data in processed folder."

$wsChangelog.Rows.Item(4).RowHeight = 28.8

# ---------------------------------------------------------------------------
# 2) column_info sheet: insert a new "Expected Format" column (new column C),
#    pushing the old "Notes" column from C to D.
# ---------------------------------------------------------------------------

$bWidth = $wsColumnInfo.Columns.Item(2).ColumnWidth

$wsColumnInfo.Columns.Item(3).Insert()
$wsColumnInfo.Columns.Item(3).ColumnWidth = $bWidth

$wsColumnInfo.Range("C1").Value = "Expected Format"

# Fill in the datetime-formatted rows first (clock in/out, both lunches)
$wsColumnInfo.Range("C4").Value = "YYYY-MM-DD HH:MM:SS"   # clock_in
$wsColumnInfo.Range("C5").Value = "YYYY-MM-DD HH:MM:SS"   # clock_out
$wsColumnInfo.Range("C6").Value = "YYYY-MM-DD HH:MM:SS"   # lunch_start
$wsColumnInfo.Range("C7").Value = "YYYY-MM-DD HH:MM:SS"   # lunch_end
$wsColumnInfo.Range("C8").Value = "YYYY-MM-DD HH:MM:SS"   # second_lunch_start
$wsColumnInfo.Range("C9").Value = "YYYY-MM-DD HH:MM:SS"   # second_lunch_end

# Then the plain-date rows
$wsColumnInfo.Range("C3").Value  = "YYYY-MM-DD"           # date
$wsColumnInfo.Range("C14").Value = "YYYY-MM-DD"           # pay_date

# Then the boolean rows
$wsColumnInfo.Range("C15").Value = "TRUE/FALSE"           # first_meal_waiver_signed
$wsColumnInfo.Range("C16").Value = "TRUE/FALSE"           # second_meal_waiver_signed
$wsColumnInfo.Range("C17").Value = "TRUE/FALSE"           # rest_break_acknowledged

# Update the (now shifted-to-D) Notes text for employment_status (row 18)
$wsColumnInfo.Range("D18").Value = "e.g., Full-Time, Part-Time, Temp, Seasonal, Contractor"

# ---------------------------------------------------------------------------
# 3) View / selection state: column_info becomes the active/selected sheet.
# ---------------------------------------------------------------------------

$wsChangelog.Range("A5").Select()

$wsColumnInfo.Activate()
$wsColumnInfo.Range("C24").Select()
